$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains two brand-new "movement" columns (move_speed / movement).
# In the real edit these were inserted just before the old "prefab" column
# (pushing prefab from H to J), and the pre-existing base_strength /
# strength_bonuses columns (which used to live in F/G) were relocated into
# the freshly inserted H/I slots, while F/G were repurposed for the new
# move_speed/movement data. Doing the insert at H (rather than at F) lets
# F and G keep their original <col> width entries (15 / 18.25) untouched.

# Step 1: insert two blank columns before the old "prefab" column (H),
# shifting prefab from H to J. F and G (and their widths) are untouched.
$ws.Range("H1:I1").EntireColumn.Insert()

# Step 2: relocate the old base_strength / strength_bonuses block (still
# sitting in F/G after the insert) into the newly inserted H/I columns.
$ws.Range("H1").Value = $ws.Range("F1").Text
$ws.Range("I1").Value = $ws.Range("G1").Text
$ws.Range("H2").Value = $ws.Range("F2").Text
$ws.Range("I2").Value = $ws.Range("G2").Text
$ws.Range("H3").Value = $ws.Range("F3").Text
$ws.Range("I3").Value = $ws.Range("G3").Text
$ws.Range("H4").Value = $ws.Range("F4").Value2
$ws.Range("I4").Value = $ws.Range("G4").Value2
$ws.Range("H5").Value = $ws.Range("F5").Value2
$ws.Range("I5").Value = $ws.Range("G5").Value2

# Step 3: write the new move_speed / movement column data into F/G.
$ws.Range("F1").Value = "move_speed"
$ws.Range("G1").Value = "movement"
$ws.Range("F2").Value = "float"
$ws.Range("G2").Value = "string"
$ws.Range("F3").Value = "基础移动速度"
$ws.Range("G3").Value = "移动方式"
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = "SimpleMove"
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = "SimpleMove"

# Step 4: give the relocated / new columns their final widths (H/I are new
# columns so there's no original width to preserve; J gets a new width too).
$ws.Columns("H").ColumnWidth = 12.12
$ws.Columns("I").ColumnWidth = 15.39
$ws.Columns("J").ColumnWidth = 36.75

# Step 5: match the saved selection / view state of the edited workbook.
$ws.Range("G10").Select()
